$wb = $excel.ActiveWorkbook

# ---- AR sheet ----
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.0002289611392375032
$ws.Range("B3").Value = 0.8271785388416928
$ws.Range("B4").Value = 5.33615690253025
$ws.Range("B5").Value = "[1.0, 0.2914812907150932, 0.26092000021777906, 0.22202665153633735, 0.2110026204074901, 0.45398322334361524, 0.23302149988310503, 0.2560351723791837, 0.18688390274123662, 0.14701871252361914, 0.18696948821120724, 0.18656123264201557, 0.18291791082120903, 0.19599752324979255, 0.1974663976038935, 0.17256564393566998, 0.16207684849081158, 0.1605062838965109, 0.17369486475797963, 0.17287800225442282]"

# ---- SETAR sheet ----
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.5557635018384144
$ws.Range("B5").Value = 0.7961118740063521
$ws.Range("B6").Value = 4.088603609273277
$ws.Range("B7").Value = 0.4941965824011549
$ws.Range("B8").Value = 0.7972885399031959
$ws.Range("B9").Value = 3.447532877435789
$ws.Range("B10").Value = "[1.0, 0.3202702857878455, 0.2712375117199265, 0.23105405910712778, 0.24627562294887398, 0.33380320512986333, 0.19516080211203352, 0.21044477236579429, 0.15697513766180268, 0.13274579997164299, 0.1559705546967117, 0.1547566796240514, 0.15393624813074946, 0.1445356178017622, 0.20585928729841005, 0.17180024540096273, 0.1517150433957696, 0.166482702206787, 0.1733166508038007, 0.16093672146961815]"

# ---- GARCH sheet ----
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.005774339871057947
$ws.Range("B3").Value = 0.0409875451045271
$ws.Range("B4").Value = 0.08838953214719292
$ws.Range("B5").Value = 0.9082342203346515
$ws.Range("B6").Value = "[1.0, 0.0030509032690374236, -0.02022154079803692, -0.038662041413097974, -0.052723473264746516, 0.3324274880435294, 0.0019122013979398704, 0.035972654613499806, -0.009691623847001118, -0.043699358679717526, -0.040289987379811605, -0.041990249187392234, -0.00961580961473966, 0.0024038732618662183, 0.02216000555938055, -0.015883014027035108, -0.030376262656352074, -0.01422940284395771, -0.030067662423678282, -0.01683654139776651]"

# ---- TARCH sheet ----
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.005552486346230828
$ws.Range("B3").Value = 0.0408483729082158
$ws.Range("B4").Value = 0.09888690436505052
$ws.Range("B5").Value = 0.9083500676710241
$ws.Range("B6").Value = "[1.0, 0.0034520740464585612, -0.020415288098801776, -0.03771021346683173, -0.05307290849185561, 0.3329332254586408, 0.001974190010471875, 0.03516175650508935, -0.009366950113707766, -0.043866800398444296, -0.04032037332794623, -0.04246427656615455, -0.009803018256398273, 0.002455552588576124, 0.02194670065395235, -0.016913548549884472, -0.031303037311896396, -0.015286673720274925, -0.03029790906115299, -0.017375639140293105]"
$ws.Range("B7").Value = -0.02090340650494344

# ---- AR_TARCH sheet ----
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.002388926979144996
$ws.Range("B3").Value = 0.03734933883024575
$ws.Range("B4").Value = 0.1094880274470461
$ws.Range("B5").Value = 0.9091213877207215
$ws.Range("B6").Value = "[1.0, 0.03996640751065263, 0.011466636625873262, 0.00802417288529159, -0.03603129864679356, 0.23302027389296048, -0.006076336095812933, 0.022577902729887647, -0.022042599724704612, -0.05032300174458967, -0.034983813365268825, -0.045108347516958955, -0.020756557136788503, -0.008586217347643942, 0.01946002372611996, -0.011390244912763157, -0.03536396712327789, -0.022984673336963978, -0.021317482107866274, -0.02011144444571112]"
$ws.Range("B7").Value = -0.04402215500150906
$ws.Range("B9").Value = 0.7929683102293414
